# Update "Results File.xlsx" with the measured throughput / transfer-ratio
# results for the MQTT QoS1 and MQTT QoS2 rows, and refresh the HTTP
# transfer-ratio figures, per "excel sheet updated" in the commit history.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- MQTT QoS1 (row 4) : Throughput (B4:I4) ----------------------------
$ws.Range("B4").Value = 83539.399999999994
$ws.Range("C4").Value = 17885.09
$ws.Range("D4").Value = 828635.17299999995
$ws.Range("E4").Value = 141518.14199999999
$ws.Range("F4").Value = 844990176.227
$ws.Range("G4").Value = 170006467.08500001
$ws.Range("H4").Value = 7041138301.8100004
$ws.Range("I4").Value = 3214330801.8800001

# ---- MQTT QoS2 (row 5) : Throughput (B5:I5) -----------------------------
$ws.Range("B5").Value = 1668.6681000000001
$ws.Range("C5").Value = 682.29
$ws.Range("D5").Value = 17726.009999999998
$ws.Range("E5").Value = 7938.2
$ws.Range("F5").Value = 46784083.530000001
$ws.Range("G5").Value = 16843543.767000001
$ws.Range("H5").Value = 511400284.29500002
$ws.Range("I5").Value = 222389616.449

# ---- HTTP (row 7) : updated transfer-ratio figures (J7:M7) -------------
$ws.Range("J7").Value = 0.8
$ws.Range("K7").Value = 81.92
$ws.Range("L7").Value = 8388.6080000000002
$ws.Range("M7").Value = 82561.296000000002

# ---- Columns B:I now hold real numbers, so Excel widened them to fit ---
$ws.Columns.Item(2).ColumnWidth = 13.998697916666666
$ws.Columns.Item(3).ColumnWidth = 13.330729166666666
$ws.Columns.Item(4).ColumnWidth = 12.998697916666666
$ws.Columns.Item(5).ColumnWidth = 13.166666666666666
$ws.Columns.Item(6).ColumnWidth = 14.498697916666666
$ws.Columns.Item(7).ColumnWidth = 13.666666666666666
$ws.Columns.Item(8).ColumnWidth = 20.330729166666668
$ws.Columns.Item(9).ColumnWidth = 21.498697916666668

# ---- Minor row height touch-up (header separators) ---------------------
$ws.Rows.Item(3).RowHeight = 16
$ws.Rows.Item(7).RowHeight = 16

# ---- Leave the selection on the first newly-entered cell ---------------
$ws.Range("H4").Select() | Out-Null
